{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// The document ends with:\n//   ... (empty paragraph)\n//   \u6d4b\u8bd5git diff\u4e2d\u5206\u522b\u6bd4\u8f83\u5de5\u4f5c\u533a\uff0c\u6682\u5b58\u533a\uff0c\u672c\u5730\u4ed3\u5e93\u7684\u547d\u4ee4\u3002  <-- last paragraph (has _GoBack bookmark)\n//\n// The edit inserts two new paragraphs right after the empty paragraph\n// (i.e. right before the final paragraph), and then replaces the text of\n// the final paragraph (keeping its formatting / bookmark) with new text.\nconst items = paragraphs.items;\nconst anchor = items[items.length - 2]; // the empty paragraph just before the last one\nconst finalParagraph = items[items.length - 1]; // keeps its pPr + bookmark\n\nconst firstNew = anchor.insertParagraph(\n  \"\u6d4b\u8bd5git diff\u4e2d\u5206\u522b\u6bd4\u8f83\u5de5\u4f5c\u533a\uff0c\u6682\u5b58\u533a\uff0c\u672c\u5730\u4ed3\u5e93\u7684\u547d\u4ee4\u3002\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nfirstNew.insertParagraph(\n  \"\u518d\u4fee\u6539\u4e00\u4e0b\uff0c\u4e4b\u524d\u7684\u4fee\u6539\u5df2\u7ecfcommit\uff0c\u672c\u4fee\u6539\u4e4b\u540e\u4ec5add\u3002\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Replace the text of the original final paragraph, preserving its\n// paragraph/run formatting and the _GoBack bookmark.\nfinalParagraph.insertText(\"\u518d\u4fee\u6539\u4e00\u4e0b\uff0c\u7136\u540e\u4ec0\u4e48\u4e5f\u4e0d\u505a\u3002\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document currently ends with:\n#   4) This is test_word2.docx file.\n#   5) (empty paragraph)\n#   6) \u6d4b\u8bd5git diff\u4e2d\u5206\u522b\u6bd4\u8f83\u5de5\u4f5c\u533a\uff0c\u6682\u5b58\u533a\uff0c\u672c\u5730\u4ed3\u5e93\u7684\u547d\u4ee4\u3002  (has the _GoBack bookmark)\n#\n# Insert two new paragraphs right after paragraph 5 (the empty one), which\n# pushes the final paragraph down while keeping its own formatting/bookmark\n# intact, then overwrite the final paragraph's text with the new wording.\n\n$anchor = $d.Paragraphs.Item(5)\n$anchor.Range.InsertParagraphAfter()\n\n$firstNew = $d.Paragraphs.Item(6)\n$firstNew.Range.Text = \"\u6d4b\u8bd5git diff\u4e2d\u5206\u522b\u6bd4\u8f83\u5de5\u4f5c\u533a\uff0c\u6682\u5b58\u533a\uff0c\u672c\u5730\u4ed3\u5e93\u7684\u547d\u4ee4\u3002\"\n\n$firstNew.Range.InsertParagraphAfter()\n$secondNew = $d.Paragraphs.Item(7)\n$secondNew.Range.Text = \"\u518d\u4fee\u6539\u4e00\u4e0b\uff0c\u4e4b\u524d\u7684\u4fee\u6539\u5df2\u7ecfcommit\uff0c\u672c\u4fee\u6539\u4e4b\u540e\u4ec5add\u3002\"\n\n# The original last paragraph (with the _GoBack bookmark) is now the 8th\n# paragraph; replace its text while keeping its paragraph/run formatting\n# and bookmark untouched.\n$count = $d.Paragraphs.Count\n$finalParagraph = $d.Paragraphs.Item($count)\n$finalParagraph.Range.Text = \"\u518d\u4fee\u6539\u4e00\u4e0b\uff0c\u7136\u540e\u4ec0\u4e48\u4e5f\u4e0d\u505a\u3002\"\n"}
